# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp footer (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 15:52"

# --- Province/city name swaps (rows keep their rank number in col A,
#     but the underlying shared-string label for the row changes) ---
$ws.Range("A26").Value = "Granada"
$ws.Range("A27").Value = "Cantabria"
$ws.Range("A35").Value = "Jaen"
$ws.Range("A36").Value = "Guadalajara"

# --- Updated numeric data (Casos totales, Casos activos, Recuperados, Muertes) ---

# Malaga
$ws.Range("B18").Value = 2143
$ws.Range("C18").Value = 523
$ws.Range("D18").Value = 1450
$ws.Range("E18").Value = 170

# Sevilla
$ws.Range("B22").Value = 1947
$ws.Range("C22").Value = 221
$ws.Range("D22").Value = 1567
$ws.Range("E22").Value = 159

# Granada (row 26 after name swap)
$ws.Range("B26").Value = 1725
$ws.Range("C26").Value = 272
$ws.Range("D26").Value = 1296
$ws.Range("E26").Value = 157

# Cantabria (row 27 after name swap)
$ws.Range("B27").Value = 1719
$ws.Range("C27").Value = 281
$ws.Range("D27").Value = 1331
$ws.Range("E27").Value = 107

# Tenerife
$ws.Range("B32").Value = 1269
$ws.Range("C32").Value = 291
$ws.Range("D32").Value = 915
$ws.Range("E32").Value = 63

# Cordoba
$ws.Range("B34").Value = 1136
$ws.Range("C34").Value = 161
$ws.Range("D34").Value = 926
$ws.Range("E34").Value = 49

# Jaen (row 35 after name swap)
$ws.Range("B35").Value = 1097
$ws.Range("C35").Value = 156
$ws.Range("D35").Value = 845
$ws.Range("E35").Value = 96

# Guadalajara (row 36 after name swap)
$ws.Range("B36").Value = 1056
$ws.Range("C36").Value = 2205
$ws.Range("D36").Value = 9768
$ws.Range("E36").Value = 139

# Cadiz
$ws.Range("B38").Value = 943
$ws.Range("C38").Value = 157
$ws.Range("D38").Value = 734
$ws.Range("E38").Value = 52

# Gran Canaria
$ws.Range("B48").Value = 479
$ws.Range("C48").Value = 116
$ws.Range("D48").Value = 336
$ws.Range("E48").Value = 27

# Almeria
$ws.Range("B51").Value = 403
$ws.Range("C51").Value = 83
$ws.Range("D51").Value = 288
$ws.Range("E51").Value = 32

# Huelva
$ws.Range("B52").Value = 318
$ws.Range("C52").Value = 57
$ws.Range("D52").Value = 239
$ws.Range("E52").Value = 22

# La Palma
$ws.Range("C56").Value = 11
$ws.Range("D56").Value = 57

# Fuerteventura
$ws.Range("C59").Value = 8
$ws.Range("D59").Value = 16
